$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells that are being updated to Text format so that
# numeric-looking strings (e.g. "247.07", "0.04577") are preserved exactly
# as text, matching the source scrape data instead of being parsed as numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.608.21'
$ws.Range("E2").Value = '  +4.02%  '
$ws.Range("D3").Value = '1.744.09'
$ws.Range("E3").Value = '  +4.27%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '247.07'
$ws.Range("E5").Value = '  +3.03%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.4812'
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("D8").Value = '0.2692'
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("D9").Value = '0.06252'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").Value = '1.740.85'
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '15.80'
$ws.Range("E12").Value = '  +6.06%  '
$ws.Range("D13").Value = '0.6184'
$ws.Range("E13").Value = '  +4.49%  '
$ws.Range("D14").Value = '4.504'
$ws.Range("E14").Value = '  +2.69%  '
$ws.Range("D15").Value = '77.57'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '26.602.91'
$ws.Range("E17").Value = '  +4.07%  '
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '0.000006894'
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("E20").Value = '  +2.34%  '
$ws.Range("D21").Value = '1.966.98'
$ws.Range("E21").Value = '  +4.22%  '
$ws.Range("D22").Value = '4.633'
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("D23").Value = '8.835'
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '5.350'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").Value = '135.89'
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("D26").Value = '15.40'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("D27").Value = '1.818'
$ws.Range("E27").Value = '  +4.85%  '
$ws.Range("D28").Value = '1.437'
$ws.Range("E28").Value = '  +3.81%  '
$ws.Range("D29").Value = '107.34'
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("D30").Value = '4.012'
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").Value = '3.753'
$ws.Range("E31").Value = '  +2.87%  '
$ws.Range("D32").Value = '0.07888'
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").Value = '0.04577'
$ws.Range("E33").Value = '  +6.85%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.6409'
$ws.Range("E35").Value = '  +5.29%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9977'
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9441'
$ws.Range("E37").Value = '  +6.01%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").Value = '112.93'
$ws.Range("E38").Value = '  +17.22%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '1.993'
$ws.Range("E39").Value = '  +6.77%  '
$ws.Range("D40").Value = '2.440'
$ws.Range("E40").Value = '  -6.10%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.763'
$ws.Range("E42").Value = '  +17.76%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.01508'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3915'
$ws.Range("E44").Value = '  +4.00%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1206'
$ws.Range("E45").Value = '  +7.62%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.683'
$ws.Range("E46").Value = '  +7.10%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05329'
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.970'
$ws.Range("E48").Value = '  +6.69%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.84'
$ws.Range("E49").Value = '  +2.94%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.270'
$ws.Range("E50").Value = '  +5.07%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3451'
$ws.Range("E51").Value = '  +3.15%  '
